# Aggiunto due colonne all'excel che potrebbero facilitare l'inserimento dei dati
# (two new "Misura" / measurement columns added to sheet1 and sheet4, shifting
# the existing header/first-data-row content one column to the right, plus a
# new "Errore" / "Errore delta x" propagated-error column & row on sheet4)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Lunghezza d'onda"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Lunghezza d'onda")

# Row 1 headers shift right by one column; new column A gets "Misura"
$ws1.Range("A1").Value = "Misura"
$ws1.Range("B1").Value = "delta x"
$ws1.Range("C1").Value = "N1"
$ws1.Range("D1").Value = "lambda"
$ws1.Range("E1").Value = "sigma lambda"

# Row 2 data shifts right by one column; new column B computes the measure/5,
# formulas updated to reference the new column positions
$ws1.Range("B2").Formula = "=A2/5"
$ws1.Range("C2").Value = 10
$ws1.Range("D2").Formula = "=2*`$B`$21*B2/C2"
$ws1.Range("E2").Formula = "=(2*`$B`$21*`$B`$22/C2)"

# ---------------------------------------------------------------------------
# Sheet 2: "Indice di rifrazione" -- no data changes, selection only
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Indice di rifrazione")

# ---------------------------------------------------------------------------
# Sheet 3: "Lunghezza pacchetti d'onda" -- no data changes, selection only
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Lunghezza pacchetti d'onda")

# ---------------------------------------------------------------------------
# Sheet 4: "Doppietto del sodio"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Doppietto del sodio")

# Row 1 headers shift right by one column; new column A gets "Misura " and a
# new column E gets "Errore"
$ws4.Range("A1").Value = "Misura "
$ws4.Range("B1").Value = "delta x"
$ws4.Range("C1").Value = "m"
$ws4.Range("D1").Value = "delta lambda"
$ws4.Range("E1").Value = "Errore"

# Row 2 data shifts right by one column; new column B computes measure/5 and
# new column E computes the propagated error
$ws4.Range("B2").Formula = "=A2/5"
$ws4.Range("D2").Formula = "=C2*B24/(2*B2)"
$ws4.Range("E2").Formula = "=D2/B2*B25"

# New row 25: propagated delta-x error, pulled from sheet 1's error constant
$ws4.Range("A25").Value = "Errore delta x"
$ws4.Range("B25").Formula = "='Lunghezza d''onda'!B22"

# ---------------------------------------------------------------------------
# Restore per-sheet selection / active-cell bookmarks to match the edit
# ---------------------------------------------------------------------------
$ws2.Range("B2").Select()
$ws3.Range("L30").Select()
$ws1.Range("D36").Select()
$ws4.Range("E3").Select()
$ws4.Activate()
